# Scheduled runner update: refresh cached market-price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve sheets.
# Each block below targets one worksheet and rewrites only the cells whose
# underlying price data changed; cells that should end up blank use
# ClearContents() instead of Value = 0 so the cell itself is removed.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1398.5454  # was 1526.5
$ws.Range("I28").Value = 268.2857  # was 293.16666
$ws.Range("K28").Value = 268.2857  # was 293.16666
$ws.Range("M28").Value = 216.7143  # was 191.83334
$ws.Range("H111").Value = 6171.7646  # was 4846.3184
$ws.Range("I111").Value = 4613.909  # was 3451.6
$ws.Range("J111").Value = 9027.833000000001  # was 7835
$ws.Range("K111").Value = 13841.727  # was 10354.8
$ws.Range("L111").Value = 27083.499  # was 23505
$ws.Range("M111").Value = -10774.727  # was -7287.799999999999
$ws.Range("N111").Value = -33217.499  # was -29639
$ws.Range("H112").Value = 5248.8237  # was 5315.222
$ws.Range("J112").Value = 5483.125  # was 5539.647
$ws.Range("L112").Value = 16449.375  # was 16618.941
$ws.Range("N112").Value = -18665.375  # was -18834.941
$ws.Range("H132").Value = 4041.111  # was 4190.5884
$ws.Range("I132").Value = 4031.4119  # was 4189.625
$ws.Range("K132").Value = 12094.2357  # was 12568.875
$ws.Range("M132").Value = -9564.235700000001  # was -10038.875
$ws.Range("H136").Value = 69744.25  # was 69991.5
$ws.Range("J136").Value = 69744.25  # was 69991.5
$ws.Range("L136").Value = 69744.25  # was 69991.5
$ws.Range("N136").Value = -79944.25  # was -80191.5
$ws.Range("H137").Value = 2447.0476  # was 2270.3333
$ws.Range("I137").Value = 1512.4286  # was 1368.7
$ws.Range("K137").Value = 4537.2858  # was 4106.1
$ws.Range("M137").Value = -1987.2858  # was -1556.1
$ws.Range("H138").Value = 4093.5642  # was 4061.8096
$ws.Range("J138").Value = 4526.4243  # was 4453.3057
$ws.Range("L138").Value = 13579.2729  # was 13359.9171
$ws.Range("N138").Value = -23859.2729  # was -23639.9171

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0  # was 22044
$ws.Range("I44").Value = 0  # was 22044
$ws.Range("K44").Value = 0  # was 22044
$ws.Range("M44").ClearContents()  # was -21556
$ws.Range("H45").Value = 11267.2  # was 10533.481
$ws.Range("I45").Value = 15557  # was 13887
$ws.Range("K45").Value = 15557  # was 13887
$ws.Range("M45").Value = -15180  # was -13510
$ws.Range("H119").Value = 57661.375  # was 61648.5
$ws.Range("J119").Value = 57661.375  # was 61648.5
$ws.Range("L119").Value = 57661.375  # was 61648.5
$ws.Range("N119").Value = -67337.375  # was -71324.5
$ws.Range("H125").Value = 0  # was 86499
$ws.Range("J125").Value = 0  # was 86499
$ws.Range("L125").Value = 0  # was 86499
$ws.Range("N125").ClearContents()  # was -96339

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 69997.25  # was 69997.75
$ws.Range("J133").Value = 69997.25  # was 69997.75
$ws.Range("L133").Value = 69997.25  # was 69997.75
$ws.Range("N133").Value = -80117.25  # was -80117.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 562.8889  # was 601.9091
$ws.Range("I22").Value = 562.8  # was 601.9091
$ws.Range("J22").Value = 563.3333  # was 0
$ws.Range("K22").Value = 562.8  # was 601.9091
$ws.Range("L22").Value = 563.3333  # was 0
$ws.Range("M22").Value = -212.8  # was -251.9091
$ws.Range("N22").Value = -1263.3333  # was absent (new cell)
$ws.Range("H52").Value = 52900  # was 53899
$ws.Range("J52").Value = 52900  # was 53899
$ws.Range("L52").Value = 52900  # was 53899
$ws.Range("N52").Value = -53488  # was -54487
$ws.Range("H122").Value = 4195  # was 3041.375
$ws.Range("I122").Value = 4397  # was 3110.7273
$ws.Range("J122").Value = 3656.3333  # was 2888.8
$ws.Range("K122").Value = 13191  # was 9332.1819
$ws.Range("L122").Value = 10968.9999  # was 8666.400000000001
$ws.Range("M122").Value = -10741  # was -6882.1819
$ws.Range("N122").Value = -15868.9999  # was -13566.4
$ws.Range("H123").Value = 89997  # was 89993.5
$ws.Range("J123").Value = 89997  # was 89993.5
$ws.Range("L123").Value = 89997  # was 89993.5
$ws.Range("N123").Value = -99797  # was -99793.5
$ws.Range("H132").Value = 3948.8794  # was 4129.5273
$ws.Range("I132").Value = 3898.311  # was 4091
$ws.Range("J132").Value = 4123.923  # was 4267.5835
$ws.Range("K132").Value = 11694.933  # was 12273
$ws.Range("L132").Value = 12371.769  # was 12802.7505
$ws.Range("M132").Value = -9164.933000000001  # was -9743
$ws.Range("N132").Value = -17431.769  # was -17862.7505
$ws.Range("H133").Value = 74323.5  # was 73764.664
$ws.Range("J133").Value = 77332.664  # was 77999
$ws.Range("L133").Value = 77332.664  # was 77999
$ws.Range("N133").Value = -82392.664  # was -83059

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1110.375  # was 1087.5883
$ws.Range("J2").Value = 536  # was 582.75
$ws.Range("L2").Value = 3216  # was 3496.5
$ws.Range("N2").Value = -3442  # was -3722.5
$ws.Range("H19").Value = 5122.5  # was 5397.5
$ws.Range("I19").Value = 930  # was 1296.6666
$ws.Range("K19").Value = 2790  # was 3889.9998
$ws.Range("M19").Value = -2616  # was -3715.9998
$ws.Range("H24").Value = 4797  # was 4047.5
$ws.Range("J24").Value = 5921.25  # was 4797
$ws.Range("L24").Value = 17763.75  # was 14391
$ws.Range("N24").Value = -18223.75  # was -14851
$ws.Range("H50").Value = 1176.4  # was 1223.5
$ws.Range("I50").Value = 1176.4  # was 1298
$ws.Range("J50").Value = 0  # was 1000
$ws.Range("K50").Value = 3529.2  # was 3894
$ws.Range("L50").Value = 0  # was 3000
$ws.Range("M50").Value = -3048.2  # was -3413
$ws.Range("N50").ClearContents()  # was -3962
$ws.Range("H53").Value = 1176.4  # was 1223.5
$ws.Range("I53").Value = 1176.4  # was 1298
$ws.Range("J53").Value = 0  # was 1000
$ws.Range("K53").Value = 3529.2  # was 3894
$ws.Range("L53").Value = 0  # was 3000
$ws.Range("M53").Value = -3048.2  # was -3413
$ws.Range("N53").ClearContents()  # was -3962
$ws.Range("H98").Value = 1351.5883  # was 1362.5294
$ws.Range("J98").Value = 1563.625  # was 1586.875
$ws.Range("L98").Value = 4690.875  # was 4760.625
$ws.Range("N98").Value = -7686.875  # was -7756.625
$ws.Range("H121").Value = 218067.22  # was 238723.95
$ws.Range("J121").Value = 715723.7  # was 1001544.6
$ws.Range("L121").Value = 2147171.1  # was 3004633.8
$ws.Range("N121").Value = -2149791.1  # was -3007253.8
$ws.Range("H137").Value = 4624.706  # was 4653.8823
$ws.Range("I137").Value = 4040.6667  # was 4248.8887
$ws.Range("J137").Value = 5281.75  # was 5109.5
$ws.Range("K137").Value = 12122.0001  # was 12746.6661
$ws.Range("L137").Value = 15845.25  # was 15328.5
$ws.Range("M137").Value = -7022.000100000001  # was -7646.666100000002
$ws.Range("N137").Value = -26045.25  # was -25528.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 434.52942  # was 487.26315
$ws.Range("I97").Value = 411.73334  # was 463
$ws.Range("J97").Value = 605.5  # was 555.2
$ws.Range("K97").Value = 411.73334  # was 463
$ws.Range("L97").Value = 605.5  # was 555.2
$ws.Range("M97").Value = 84.26666  # was 33
$ws.Range("N97").Value = -1597.5  # was -1547.2
$ws.Range("H99").Value = 16601.25  # was 19690
$ws.Range("I99").Value = 8973  # was 9628.200000000001
$ws.Range("K99").Value = 8973  # was 9628.200000000001
$ws.Range("M99").Value = -6727  # was -7382.200000000001
$ws.Range("H102").Value = 2412.25  # was 2593.1428
$ws.Range("I102").Value = 2412.25  # was 2593.1428
$ws.Range("K102").Value = 2412.25  # was 2593.1428
$ws.Range("M102").Value = -790.25  # was -971.1428000000001
$ws.Range("H122").Value = 3743.65  # was 3954.2222
$ws.Range("I122").Value = 3261.625  # was 3359.2666
$ws.Range("J122").Value = 5671.75  # was 6929
$ws.Range("K122").Value = 9784.875  # was 10077.7998
$ws.Range("L122").Value = 17015.25  # was 20787
$ws.Range("M122").Value = -7334.875  # was -7627.799800000001
$ws.Range("N122").Value = -21915.25  # was -25687

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10264.88  # was 10507.5
$ws.Range("I40").Value = 7704.3887  # was 7720.5557
$ws.Range("J40").Value = 16849  # was 18868.334
$ws.Range("K40").Value = 7704.3887  # was 7720.5557
$ws.Range("L40").Value = 16849  # was 18868.334
$ws.Range("M40").Value = -7568.3887  # was -7584.5557
$ws.Range("N40").Value = -17121  # was -19140.334
$ws.Range("H100").Value = 5083.3335  # was 4571.357
$ws.Range("I100").Value = 3250  # was 2666.5
$ws.Range("K100").Value = 3250  # was 2666.5
$ws.Range("M100").Value = -2709  # was -2125.5
$ws.Range("H132").Value = 2902.2686  # was 2689.2534
$ws.Range("I132").Value = 2932.451  # was 2724.8276
$ws.Range("J132").Value = 2806.0625  # was 2530.5386
$ws.Range("K132").Value = 8797.352999999999  # was 8174.4828
$ws.Range("L132").Value = 8418.1875  # was 7591.6158
$ws.Range("M132").Value = -6267.352999999999  # was -5644.4828
$ws.Range("N132").Value = -13478.1875  # was -12651.6158
$ws.Range("H134").Value = 85999  # was 81499.5
$ws.Range("J134").Value = 85999  # was 81499.5
$ws.Range("L134").Value = 85999  # was 81499.5
$ws.Range("N134").Value = -96139  # was -91639.5
$ws.Range("H141").Value = 76255.39999999999  # was 76275.2
$ws.Range("J141").Value = 76255.39999999999  # was 76275.2
$ws.Range("L141").Value = 76255.39999999999  # was 76275.2
$ws.Range("N141").Value = -86615.39999999999  # was -86635.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4325.0557  # was 4325.1665
$ws.Range("I107").Value = 3966.25  # was 3761.2307
$ws.Range("J107").Value = 5042.6665  # was 5791.4
$ws.Range("K107").Value = 11898.75  # was 11283.6921
$ws.Range("L107").Value = 15127.9995  # was 17374.2
$ws.Range("M107").Value = -9978.75  # was -9363.6921
$ws.Range("N107").Value = -18967.9995  # was -21214.2
$ws.Range("H135").Value = 64588  # was 62967.223
$ws.Range("J135").Value = 64588  # was 62967.223
$ws.Range("L135").Value = 64588  # was 62967.223
$ws.Range("N135").Value = -74728  # was -73107.223
